$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format before writing, so numeric-looking
# price strings (e.g. "1.000", "241.32") are stored as text, matching
# the workbook's existing inlineStr cells instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '29.051.31'
$ws.Range("E2").Value = '  -0.47%  '
# Row 3
$ws.Range("D3").Value = '1.825.94'
$ws.Range("E3").Value = '  -0.54%  '
# Row 4
$ws.Range("D4").Value = '0.9993'
$ws.Range("E4").Value = '  -0.02%  '
# Row 5
$ws.Range("D5").Value = '241.32'
$ws.Range("E5").Value = '  +0.17%  '
# Row 6
$ws.Range("D6").Value = '0.6357'
$ws.Range("E6").Value = '  -4.48%  '
# Row 7
$ws.Range("E7").Value = '  +0.01%  '
# Row 8
$ws.Range("D8").Value = '44.78'
$ws.Range("E8").Value = '  +6.70%  '
# Row 9
$ws.Range("D9").Value = '0.2935'
$ws.Range("E9").Value = '  +0.51%  '
# Row 10
$ws.Range("D10").Value = '0.07337'
$ws.Range("E10").Value = '  -0.37%  '
# Row 11
$ws.Range("D11").Value = '22.80'
$ws.Range("E11").Value = '  +0.76%  '
# Row 12
$ws.Range("D12").Value = '0.07668'
$ws.Range("E12").Value = '  -0.64%  '
# Row 13
$ws.Range("D13").Value = '1.827.28'
$ws.Range("E13").Value = '  -0.28%  '
# Row 14
$ws.Range("D14").Value = '4.986'
$ws.Range("E14").Value = '  +0.09%  '
# Row 15
$ws.Range("D15").Value = '0.6634'
$ws.Range("E15").Value = '  -0.66%  '
# Row 16
$ws.Range("E16").Value = '  -1.29%  '
# Row 17
$ws.Range("D17").Value = '0.000008662'
$ws.Range("E17").Value = '  +4.92%  '
# Row 18
$ws.Range("D18").Value = '6.032'
$ws.Range("E18").Value = '  -1.40%  '
# Row 19
$ws.Range("D19").Value = '28.900.76'
$ws.Range("E19").Value = '  -0.76%  '
# Row 20
$ws.Range("D20").Value = '2.075.23'
$ws.Range("E20").Value = '  +1.19%  '
# Row 21
$ws.Range("D21").Value = '224.51'
$ws.Range("E21").Value = '  -0.43%  '
# Row 22
$ws.Range("D22").Value = '12.38'
$ws.Range("E22").Value = '  -0.50%  '
# Row 23
$ws.Range("E23").Value = '  -0.03%  '
# Row 24
$ws.Range("D24").Value = '7.113'
$ws.Range("E24").Value = '  -0.05%  '
# Row 25
$ws.Range("D25").Value = '1.000'
$ws.Range("E25").Value = '  -0.02%  '
# Row 26
$ws.Range("D26").Value = '158.04'
$ws.Range("E26").Value = '  -1.62%  '
# Row 27
$ws.Range("D27").Value = '8.465'
$ws.Range("E27").Value = '  -1.82%  '
# Row 28
$ws.Range("D28").Value = '0.1371'
$ws.Range("E28").Value = '  -1.49%  '
# Row 29
$ws.Range("D29").Value = '17.92'
$ws.Range("E29").Value = '  -0.08%  '
# Row 30
$ws.Range("D30").Value = '1.502'
$ws.Range("E30").Value = '  -0.75%  '
# Row 31
$ws.Range("D31").Value = '4.089'
$ws.Range("E31").Value = '  -0.39%  '
# Row 32
$ws.Range("D32").Value = '4.027'
$ws.Range("E32").Value = '  -0.06%  '
# Row 33
$ws.Range("E33").Value = '  +1.84%  '
# Row 34
$ws.Range("D34").Value = '0.05291'
$ws.Range("E34").Value = '  -0.16%  '
# Row 35
$ws.Range("D35").Value = '1.838'
$ws.Range("E35").Value = '  -1.72%  '
# Row 36
$ws.Range("D36").Value = '0.7360'
$ws.Range("E36").Value = '  -2.18%  '
# Row 37
$ws.Range("E37").Value = '  +2.19%  '
# Row 38
$ws.Range("D38").Value = '2.649'
$ws.Range("E38").Value = '  -1.08%  '
# Row 39
$ws.Range("D39").Value = '1.295.78'
$ws.Range("E39").Value = '  +0.15%  '
# Row 40
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = '0.01782'
$ws.Range("E40").Value = '  -0.66%  '
# Row 41
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").Value = '2.741'
$ws.Range("E41").Value = '  +0.78%  '
# Row 42
$ws.Range("D42").Value = '6.300'
$ws.Range("E42").Value = '  +5.80%  '
# Row 43
$ws.Range("D43").Value = '0.8995'
$ws.Range("E43").Value = '  -2.26%  '
# Row 44
$ws.Range("D44").Value = '0.9996'
$ws.Range("E44").Value = '  -0.62%  '
# Row 45
$ws.Range("D45").Value = '102.49'
$ws.Range("E45").Value = '  +0.28%  '
# Row 46
$ws.Range("D46").Value = '1.974.22'
$ws.Range("E46").Value = '  +0.75%  '
# Row 47
$ws.Range("E47").Value = '  -0.50%  '
# Row 48
$ws.Range("D48").Value = '63.99'
$ws.Range("E48").Value = '  +1.37%  '
# Row 49
$ws.Range("D49").Value = '0.00000000120'
$ws.Range("E49").Value = '  -0.78%  '
# Row 50
$ws.Range("E50").Value = '  -2.15%  '
# Row 51
$ws.Range("B51").Value = 'XinFinNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range("D51").Value = '0.07264'
$ws.Range("E51").Value = '  -15.75%  '

# Restore default (unstyled) appearance for column D now that the
# text number-format has served its purpose of preventing numeric coercion.
$ws.Range("D2:D51").Style = "Normal"
